$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the sample-data header row (ID_REF + GSM2299032..GSM2299043) — this
# also shifts the rows below it (39 -> 38, 42 -> 41) and renumbers the
# shared-string indices used by B20:D20 automatically.
$ws.Rows.Item(35).Delete()

# Update the active selection to match the saved view state.
$ws.Range("B29").Select()
